$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 586
$ws1.Range("F9").Value = 201
$ws1.Range("F13").Value = 135
$ws1.Range("F16").Value = 5209
$ws1.Range("F17").Value = 105
$ws1.Range("F21").Value = 274
$ws1.Range("F23").Value = 6166
$ws1.Range("F28").Value = 14918
$ws1.Range("F33").Value = 10873
$ws1.Range("F34").Value = 690
$ws1.Range("F35").Value = 4251

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 586
$ws4.Range("F9").Value = 201
$ws4.Range("F13").Value = 135
$ws4.Range("F17").Value = 5209
$ws4.Range("F18").Value = 105
$ws4.Range("F23").Value = 274
$ws4.Range("F24").Value = 37
$ws4.Range("F26").Value = 6166
$ws4.Range("F31").Value = 14918
$ws4.Range("F36").Value = 10873
$ws4.Range("F37").Value = 690
$ws4.Range("F38").Value = 4251
